# Apply updated TPM-derived values to Plau-St14 sheet (Young D7 dataset).
# Ligand/receptor/edge statistics (columns G-J, M-P, Q-T) were recomputed
# with new TPM data; columns A-F and K-L are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.44975366666667
$ws.Range("H2").Value = 46.349261
$ws.Range("I2").Value = 0.1287486886000874
$ws.Range("J2").Value = 0.1287486886000874
$ws.Range("M2").Value = 1.483864
$ws.Range("N2").Value = 4.451592
$ws.Range("O2").Value = 0.4831229533629666
$ws.Range("P2").Value = 0.4831229533629666
$ws.Range("Q2").Value = 22.92533327483466
$ws.Range("R2").Value = 206.327999473512
$ws.Range("S2").Value = 0.06220144667808314
$ws.Range("T2").Value = 0.06220144667808313
$ws.Range("G3").Value = 15.44975366666667
$ws.Range("H3").Value = 46.349261
$ws.Range("I3").Value = 0.1287486886000874
$ws.Range("J3").Value = 0.1287486886000874
$ws.Range("O3").Value = 0.2034038545501667
$ws.Range("P3").Value = 0.2034038545501667
$ws.Range("Q3").Value = 9.651996707027111
$ws.Range("R3").Value = 86.867970363244
$ws.Range("S3").Value = 0.02618797952953688
$ws.Range("T3").Value = 0.02618797952953688
$ws.Range("G4").Value = 15.44975366666667
$ws.Range("H4").Value = 46.349261
$ws.Range("I4").Value = 0.1287486886000874
$ws.Range("J4").Value = 0.1287486886000874
$ws.Range("M4").Value = 0.6102786666666665
$ws.Range("N4").Value = 1.830836
$ws.Range("O4").Value = 0.1986972066270314
$ws.Range("P4").Value = 0.1986972066270314
$ws.Range("Q4").Value = 9.428655068021776
$ws.Range("R4").Value = 84.85789561219599
$ws.Range("S4").Value = 0.0255820047817309
$ws.Range("T4").Value = 0.02558200478173089
$ws.Range("G5").Value = 15.44975366666667
$ws.Range("H5").Value = 46.349261
$ws.Range("I5").Value = 0.1287486886000874
$ws.Range("J5").Value = 0.1287486886000874
$ws.Range("M5").Value = 0.352523
$ws.Range("N5").Value = 1.057569
$ws.Range("O5").Value = 0.1147759854598353
$ws.Range("P5").Value = 0.1147759854598353
$ws.Range("Q5").Value = 5.446393511834333
$ws.Range("R5").Value = 49.017541606509
$ws.Range("S5").Value = 0.0147772576107365
$ws.Range("T5").Value = 0.0147772576107365
$ws.Range("I6").Value = 0.4074352211478151
$ws.Range("J6").Value = 0.4074352211478151
$ws.Range("M6").Value = 1.483864
$ws.Range("N6").Value = 4.451592
$ws.Range("O6").Value = 0.4831229533629666
$ws.Range("P6").Value = 0.4831229533629666
$ws.Range("Q6").Value = 72.54899707547999
$ws.Range("R6").Value = 652.94097367932
$ws.Range("S6").Value = 0.1968413073450258
$ws.Range("T6").Value = 0.1968413073450258
$ws.Range("I7").Value = 0.4074352211478151
$ws.Range("J7").Value = 0.4074352211478151
$ws.Range("O7").Value = 0.2034038545501667
$ws.Range("P7").Value = 0.2034038545501667
$ws.Range("S7").Value = 0.08287389446096517
$ws.Range("T7").Value = 0.08287389446096517
$ws.Range("I8").Value = 0.4074352211478151
$ws.Range("J8").Value = 0.4074352211478151
$ws.Range("M8").Value = 0.6102786666666665
$ws.Range("N8").Value = 1.830836
$ws.Range("O8").Value = 0.1986972066270314
$ws.Range("P8").Value = 0.1986972066270314
$ws.Range("Q8").Value = 29.83771100533999
$ws.Range("R8").Value = 268.53939904806
$ws.Range("S8").Value = 0.08095624032353767
$ws.Range("T8").Value = 0.08095624032353767
$ws.Range("I9").Value = 0.4074352211478151
$ws.Range("J9").Value = 0.4074352211478151
$ws.Range("M9").Value = 0.352523
$ws.Range("N9").Value = 1.057569
$ws.Range("O9").Value = 0.1147759854598353
$ws.Range("P9").Value = 0.1147759854598353
$ws.Range("Q9").Value = 17.235535127235
$ws.Range("R9").Value = 155.119816145115
$ws.Range("S9").Value = 0.04676377901828641
$ws.Range("T9").Value = 0.04676377901828641
$ws.Range("G10").Value = 24.32144666666666
$ws.Range("H10").Value = 72.96433999999999
$ws.Range("I10").Value = 0.2026798893205849
$ws.Range("J10").Value = 0.2026798893205849
$ws.Range("M10").Value = 1.483864
$ws.Range("N10").Value = 4.451592
$ws.Range("O10").Value = 0.4831229533629666
$ws.Range("P10").Value = 0.4831229533629666
$ws.Range("Q10").Value = 36.08971913658666
$ws.Range("R10").Value = 324.8074722292799
$ws.Range("S10").Value = 0.09791930671584015
$ws.Range("T10").Value = 0.09791930671584015
$ws.Range("G11").Value = 24.32144666666666
$ws.Range("H11").Value = 72.96433999999999
$ws.Range("I11").Value = 0.2026798893205849
$ws.Range("J11").Value = 0.2026798893205849
$ws.Range("O11").Value = 0.2034038545501667
$ws.Range("P11").Value = 0.2034038545501667
$ws.Range("Q11").Value = 15.19445087615111
$ws.Range("R11").Value = 136.75005788536
$ws.Range("S11").Value = 0.04122587072760812
$ws.Range("T11").Value = 0.04122587072760812
$ws.Range("G12").Value = 24.32144666666666
$ws.Range("H12").Value = 72.96433999999999
$ws.Range("I12").Value = 0.2026798893205849
$ws.Range("J12").Value = 0.2026798893205849
$ws.Range("M12").Value = 0.6102786666666665
$ws.Range("N12").Value = 1.830836
$ws.Range("O12").Value = 0.1986972066270314
$ws.Range("P12").Value = 0.1986972066270314
$ws.Range("Q12").Value = 14.84286004313777
$ws.Range("R12").Value = 133.58574038824
$ws.Range("S12").Value = 0.04027192784747612
$ws.Range("T12").Value = 0.04027192784747612
$ws.Range("G13").Value = 24.32144666666666
$ws.Range("H13").Value = 72.96433999999999
$ws.Range("I13").Value = 0.2026798893205849
$ws.Range("J13").Value = 0.2026798893205849
$ws.Range("M13").Value = 0.352523
$ws.Range("N13").Value = 1.057569
$ws.Range("O13").Value = 0.1147759854598353
$ws.Range("P13").Value = 0.1147759854598353
$ws.Range("Q13").Value = 8.573869343273332
$ws.Range("R13").Value = 77.16482408946
$ws.Range("S13").Value = 0.02326278402966048
$ws.Range("T13").Value = 0.02326278402966048
$ws.Range("G14").Value = 31.33616366666666
$ws.Range("H14").Value = 94.00849099999999
$ws.Range("I14").Value = 0.2611362009315126
$ws.Range("J14").Value = 0.2611362009315126
$ws.Range("M14").Value = 1.483864
$ws.Range("N14").Value = 4.451592
$ws.Range("O14").Value = 0.4831229533629666
$ws.Range("P14").Value = 0.4831229533629666
$ws.Range("Q14").Value = 46.49860516307466
$ws.Range("R14").Value = 418.4874464676719
$ws.Range("S14").Value = 0.1261608926240174
$ws.Range("T14").Value = 0.1261608926240174
$ws.Range("G15").Value = 31.33616366666666
$ws.Range("H15").Value = 94.00849099999999
$ws.Range("I15").Value = 0.2611362009315126
$ws.Range("J15").Value = 0.2611362009315126
$ws.Range("O15").Value = 0.2034038545501667
$ws.Range("P15").Value = 0.2034038545501667
$ws.Range("Q15").Value = 19.57678776290711
$ws.Range("R15").Value = 176.191089866164
$ws.Range("S15").Value = 0.05311610983205648
$ws.Range("T15").Value = 0.05311610983205648
$ws.Range("G16").Value = 31.33616366666666
$ws.Range("H16").Value = 94.00849099999999
$ws.Range("I16").Value = 0.2611362009315126
$ws.Range("J16").Value = 0.2611362009315126
$ws.Range("M16").Value = 0.6102786666666665
$ws.Range("N16").Value = 1.830836
$ws.Range("O16").Value = 0.1986972066270314
$ws.Range("P16").Value = 0.1986972066270314
$ws.Range("Q16").Value = 19.12379218094177
$ws.Range("R16").Value = 172.114129628476
$ws.Range("S16").Value = 0.05188703367428676
$ws.Range("T16").Value = 0.05188703367428676
$ws.Range("G17").Value = 31.33616366666666
$ws.Range("H17").Value = 94.00849099999999
$ws.Range("I17").Value = 0.2611362009315126
$ws.Range("J17").Value = 0.2611362009315126
$ws.Range("M17").Value = 0.352523
$ws.Range("N17").Value = 1.057569
$ws.Range("O17").Value = 0.1147759854598353
$ws.Range("P17").Value = 0.1147759854598353
$ws.Range("Q17").Value = 11.04671842426433
$ws.Range("R17").Value = 99.42046581837899
$ws.Range("S17").Value = 0.02997216480115192
$ws.Range("T17").Value = 0.02997216480115192
